# Add season-record columns (Wins / Losses / Ties) to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): AD1="Wins", AE1="Losses", AF1="Ties" ---
# Copy the formatting of the existing last header cell (AC1) onto the
# three new header cells so they pick up the same bold/centered/bordered
# style used by the rest of row 1, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-52): every player gets the team's season record ---
$wins = 69
$losses = 93
$ties = 0

for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-52"
